$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns F and G (header + data rows 1-16):
# F header "Rank" <-> G header "Rata Rata"
# F/G data values swap per row as well.
for ($r = 1; $r -le 16; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 6).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $fVal
}
